# HubTalentos-Avaliacao / AdvantageOnline_Data.xlsx
# "Classe renomeada apenas, sem alteracao no codigo."
#
# On the "Users" sheet, the four test-account login names (column B, rows 2-5)
# are each bumped to the "next" generated value:
#   BRUNO109   -> BRUNO114
#   Marcela103 -> Marcela104
#   BrUn94     -> BrUn95
#   Maria14    -> Maria15
#
# The cells are written in the order B3, B4, B5, B2 so that the shared-string
# table ends up laid out the same way Excel produces it (new strings are
# appended in the order they are first typed, and the four strings that
# become unused get reclaimed/compacted on save).
#
# The workbook is also left with the "Users" sheet active (instead of
# "Headphone") and cell B2 selected on it.

$wb = $excel.ActiveWorkbook
$wsUsers = $wb.Worksheets.Item("Users")

$wsUsers.Cells.Item(3, 2).Value = "Marcela104"
$wsUsers.Cells.Item(4, 2).Value = "BrUn95"
$wsUsers.Cells.Item(5, 2).Value = "Maria15"
$wsUsers.Cells.Item(2, 2).Value = "BRUNO114"

$wsUsers.Activate()
$wsUsers.Range("B2").Select()
